$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trade_Param")

# The two data rows (FLO_DELIV / FLO_COST definitions) are wiped out,
# leaving only the header rows (2-4) and the formatting/styles in place.
$ws.Range("B6:F7").ClearContents()

# Row 6 no longer needs the extra height that was required to show the
# wrapped text, so let it fall back to the sheet's default height.
$ws.Rows(6).AutoFit()

# Row 7 keeps a slightly taller, thick-bottom-bordered row (matching the
# look of row 4 above it) even though its contents are now empty.
$ws.Rows(7).RowHeight = 15

# Update the active selection to reflect where the user left off editing.
$null = $ws.Range("G7").Select()

Write-Host "Cleared Trade_Param data rows (B6:F7) and updated selection."
